$d = $word.ActiveDocument

$pairs = @(
    @("19÷9=2, 1", "70÷3=23, 1"),
    @("20÷5=4, 0", "73÷3=24, 1"),
    @("62÷9=6, 8", "21÷5=4, 1"),
    @("97÷7=13, 6", "58÷4=14, 2"),
    @("24÷8=3, 0", "16÷9=1, 7"),
    @("70÷2=35, 0", "98÷6=16, 2"),
    @("82÷2=41, 0", "70÷8=8, 6"),
    @("16÷8=2, 0", "92÷6=15, 2"),
    @("61÷9=6, 7", "16÷9=1, 7"),
    @("85÷7=12, 1", "82÷4=20, 2"),
    @("79÷5=15, 4", "40÷7=5, 5"),
    @("16÷5=3, 1", "41÷7=5, 6"),
    @("30÷3=10, 0", "30÷8=3, 6"),
    @("78÷8=9, 6", "48÷6=8, 0"),
    @("97÷5=19, 2", "10÷2=5, 0"),
    @("80÷9=8, 8", "51÷2=25, 1"),
    @("93÷6=15, 3", "35÷7=5, 0"),
    @("51÷3=17, 0", "11÷9=1, 2"),
    @("91÷6=15, 1", "27÷3=9, 0"),
    @("42÷4=10, 2", "54÷2=27, 0"),
    @("98÷4=24, 2", "11÷9=1, 2"),
    @("34÷6=5, 4", "22÷8=2, 6"),
    @("73÷8=9, 1", "23÷2=11, 1"),
    @("63÷3=21, 0", "88÷7=12, 4"),
    @("81÷8=10, 1", "75÷4=18, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
